# Updated queries for C3DC first half testcases.
#
# The SQL queries stored in column B (rows 2-7) and C2 joined tables using a
# generic "id" column (e.g. std.id = prt."study.id"). They are updated to use
# the fully-qualified, entity-specific id columns (e.g.
# std.study_id = prt."study.study_id").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold one of the SQL queries touched by this change.
$targetCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

# Ordered list of (old join condition -> new join condition) replacements
# applied to every query above.
$joinReplacements = @(
    @{ Old = 'df_participant prt ON std.id = prt."study.id"';            New = 'df_participant prt ON std.study_id = prt."study.study_id"' },
    @{ Old = 'df_diagnoses dgn ON prt.id = dgn."participant.id"';        New = 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"' },
    @{ Old = 'df_treatments trt ON prt.id = trt."participant.id"';       New = 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"' },
    @{ Old = 'df_treatment_resp trr ON prt.id = trr."participant.id"';   New = 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"' },
    @{ Old = 'df_survival srv ON prt.id = srv."participant.id"';         New = 'df_survival srv ON prt.participant_id = srv."participant.participant_id"' },
    @{ Old = 'df_reference_files rfs ON std.id = rfs."study.id"';        New = 'df_reference_files rfs ON std.study_id = rfs."study.study_id"' }
)

# A few of the queries (the ones immediately followed by a WHERE clause with
# no blank/condition line in between) also lost the trailing space after
# "WHERE" on the line right after the reference_files join once re-typed.
$whereSpaceFixCells = @("C2", "B5", "B7")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2

    foreach ($rep in $joinReplacements) {
        $text = $text -replace [regex]::Escape($rep.Old), $rep.New
    }

    if ($whereSpaceFixCells -contains $addr) {
        $text = $text -replace [regex]::Escape('df_reference_files rfs ON std.study_id = rfs."study.study_id"' + "`nWHERE `n"), ('df_reference_files rfs ON std.study_id = rfs."study.study_id"' + "`nWHERE`n")
    }

    $cell.Value2 = $text
}
